# Add a new worksheet "monthly_ph" at the end of the workbook (after the
# current last sheet, "annual_temp") summarizing Philippine monthly
# climatology: min/mean/max temperature and precipitation per month.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "monthly_ph"

# Header row - set the columns that reuse existing shared strings first
# (min/mean/max already exist in the workbook), then the two brand-new
# strings (precipitation, category) so they land in the same shared-string
# order as the source edit.
$ws.Range("B1").Value = "min"
$ws.Range("C1").Value = "mean"
$ws.Range("D1").Value = "max"
$ws.Range("E1").Value = "precipitation"
$ws.Range("A1").Value = "category"

# Monthly data: category, min, mean, max, precipitation
$ws.Range("A2").Value = "Jan"
$ws.Range("B2").Value = 20.86
$ws.Range("C2").Value = 24.96
$ws.Range("D2").Value = 29.11
$ws.Range("E2").Value = 160.82

$ws.Range("A3").Value = "Feb"
$ws.Range("B3").Value = 20.83
$ws.Range("C3").Value = 25.18
$ws.Range("D3").Value = 29.59
$ws.Range("E3").Value = 118.73

$ws.Range("A4").Value = "Mar"
$ws.Range("B4").Value = 21.42
$ws.Range("C4").Value = 26.09
$ws.Range("D4").Value = 30.81
$ws.Range("E4").Value = 117.2

$ws.Range("A5").Value = "Apr"
$ws.Range("B5").Value = 22.33
$ws.Range("C5").Value = 27.09
$ws.Range("D5").Value = 31.89
$ws.Range("E5").Value = 116.61

$ws.Range("A6").Value = "May"
$ws.Range("B6").Value = 22.87
$ws.Range("C6").Value = 27.42
$ws.Range("D6").Value = 32.02
$ws.Range("E6").Value = 203.2

$ws.Range("A7").Value = "Jun"
$ws.Range("B7").Value = 22.59
$ws.Range("C7").Value = 26.85
$ws.Range("D7").Value = 31.16
$ws.Range("E7").Value = 251.66

$ws.Range("A8").Value = "Jul"
$ws.Range("B8").Value = 22.27
$ws.Range("C8").Value = 26.31
$ws.Range("D8").Value = 30.4
$ws.Range("E8").Value = 290.62

$ws.Range("A9").Value = "Aug"
$ws.Range("B9").Value = 22.3
$ws.Range("C9").Value = 26.28
$ws.Range("D9").Value = 30.31
$ws.Range("E9").Value = 314.01

$ws.Range("A10").Value = "Sep"
$ws.Range("B10").Value = 22.18
$ws.Range("C10").Value = 26.26
$ws.Range("D10").Value = 30.38
$ws.Range("E10").Value = 276.13

$ws.Range("A11").Value = "Oct"
$ws.Range("B11").Value = 22.01
$ws.Range("C11").Value = 26.17
$ws.Range("D11").Value = 30.39
$ws.Range("E11").Value = 280.33

$ws.Range("A12").Value = "Nov"
$ws.Range("B12").Value = 21.87
$ws.Range("C12").Value = 25.95
$ws.Range("D12").Value = 30.07
$ws.Range("E12").Value = 251.71

$ws.Range("A13").Value = "Dec"
$ws.Range("B13").Value = 21.49
$ws.Range("C13").Value = 25.43
$ws.Range("D13").Value = 29.42
$ws.Range("E13").Value = 254.52

# Match the "best fit" column widths the original author's Excel computed
# for the min/mean columns (B, C) once the sheet's data/headers were in
# place. (The new sheet becomes the active/selected tab automatically,
# mirroring activeTab moving from 3 to 4.)
$ws.Columns.Item(2).ColumnWidth = 15.96
$ws.Columns.Item(3).ColumnWidth = 17.67
